# Entrevista.docx - "Atualizacao de cores e testo" edit
#
# The title run "Entrevista" is replaced by the red (FF0000) text
# "Check list guia para entrevista". Word's spell-checker would have
# flagged the two English words ("Check" and "list") as not recognised
# by the pt-BR dictionary, so they get wrapped in w:proofErr spellStart/
# spellEnd markers. The pre-existing (empty) "_GoBack" bookmark, which
# originally sat at the very end of "Entrevista", ends up back in the
# same relative spot in the document - i.e. right after "...guia pa",
# before "ra entrevista" - once the new text is typed over it.

$d = $word.ActiveDocument

# Locate "Entrevista" so we don't have to hard-code character offsets.
$findRng = $d.Content
$found = $findRng.Find.Execute("Entrevista", $false, $false, $false, $false, `
                                $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Host "ERROR: 'Entrevista' not found"
}
$baseStart = $findRng.Start
$target = $d.Range($baseStart, $findRng.End)

# Common run formatting for the new (red) text.
$rPr = "<w:rPr>" + `
         "<w:rFonts w:ascii='Times New Roman' w:eastAsia='Times New Roman' " + `
           "w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" + `
         "<w:color w:val='FF0000'/>" + `
         "<w:sz w:val='24'/>" + `
         "<w:szCs w:val='24'/>" + `
         "<w:lang w:eastAsia='pt-BR'/>" + `
       "</w:rPr>"

$body = "<w:proofErr w:type='spellStart'/>" + `
        "<w:r>$rPr<w:t>Check</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r>$rPr<w:t xml:space=`"preserve`"> </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r>$rPr<w:t>list</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r>$rPr<w:t xml:space=`"preserve`"> guia pa</w:t></w:r>" + `
        "<w:r>$rPr<w:t>ra entrevista</w:t></w:r>"

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$xml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" + `
         "<pkg:part pkg:name='/word/document.xml' " + `
           "pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" + `
           "<pkg:xmlData><w:document $wns><w:body><w:p>$body</w:p></w:body></w:document></pkg:xmlData>" + `
         "</pkg:part>" + `
       "</pkg:package>"

$target.InsertXML($xml)

# Put the "_GoBack" bookmark back where it belongs: right after
# "Check list guia pa", before "ra entrevista".
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$bmPos = $baseStart + "Check list guia pa".Length
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
